$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 469; this pushes the existing rows 469-557
# down to 470-558 and keeps all of their data/formatting intact.
$ws.Rows.Item(469).Insert()

# Populate the newly inserted row 469 with the new data record.
$ws.Range("A469").Value = 6
$ws.Range("B469").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C469").Value = "Metropolitana"
$ws.Range("D469").Value = 44951
$ws.Range("E469").Value = 13
$ws.Range("F469").Value = 100112043
$ws.Range("G469").Value = "Pepino ensalada"
$ws.Range("H469").Value = "Sin especificar"
$ws.Range("I469").Value = "Primera"
$ws.Range("J469").Value = 1200
$ws.Range("K469").Value = 6000
$ws.Range("L469").Value = 7000
$ws.Range("M469").Value = 6458
$ws.Range("N469").Value = "`$/caja 60 unidades"
$ws.Range("O469").Value = "Región Metropolitana"
$ws.Range("P469").Value = 108
$ws.Range("Q469").Value = 60
$ws.Range("R469").Value = "Hortaliza"
